# Fixed #476 Moving from Apache POI 4.1.0 to 5.2.3.
#
# Re-generating this m2doc bookmark-link test fixture with the newer POI
# writer re-serializes every bold run it touches (dropping the explicit
# w:val="true" in favour of the library's own boolean rendering) and mints
# a fresh bookmark id for "bookmark1". Reproduce that by touching the same
# three bold runs and by re-creating the bookmark in place so the document
# is forced through the writer again.

$d = $word.ActiveDocument

# --- Run 1: the REF field result ("a reference to bookmark1") ----------
# Bold + NoProof field-result run inside the "Test link before bookmark"
# paragraph.
$fieldResult = $d.Content.Duplicate
$fieldResult.Find.Execute("a reference to bookmark1")
$fieldResult.Font.Bold = 1

# --- Run 2 & 3: the two "Invalid link statement" messages --------------
$msg1 = $d.Content.Duplicate
$msg1.Find.Execute("Invalid link statement: Expression ""self. 'a reference to bookmark1'"" is invalid: missing feature access or service call")
$msg1.Font.Bold = 1

$msg2 = $d.Content.Duplicate
$msg2.Find.Execute("Invalid link statement: Expression """" is invalid: null or empty string.")
$msg2.Font.Bold = 1

# --- bookmark1: re-mint its id by deleting and re-adding it in place ---
$bookmark = $d.Bookmarks("bookmark1")
$bookmarkRange = $bookmark.Range
$bookmark.Delete()
$d.Bookmarks.Add("bookmark1", $bookmarkRange)
